$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The sheet has two parts tables:
#   Table1 ("Quadcopter:")  originally A2:E16  (header row2, data rows3-15, totals row16)
#   Table2 ("Receiver:")    originally A19:E27 (header row19, data rows20-26, totals row27)
#
# Three new parts are being added to each table, in the same relative
# positions: right after the "Microprocessor" row a new "IC Socket" row
# is inserted, and right after the "10k Resistor" row two new rows
# ("3mm LED" and "100R Resistor") are inserted.
#
# We insert the blank rows from the bottom of the sheet upwards so that
# earlier row numbers used for the Insert calls stay valid.
# ------------------------------------------------------------------

# ---- Table2 (rows are below Table1, so handle first) ----
# two new rows after the "10k Resistor" row (row 24)
$ws.Rows("25:26").Insert()
# one new row after the "Microprocessor" row (row 20)
$ws.Rows("21:21").Insert()

# ---- Table1 ----
# two new rows after the "10k Resistor (x5)" row (row 6)
$ws.Rows("7:8").Insert()
# one new row after the "Microprocessor" row (row 3)
$ws.Rows("4:4").Insert()

# ------------------------------------------------------------------
# Fill in the new rows for Table1
# ------------------------------------------------------------------

# Row 4: IC Socket
$ws.Range("A4").Value = "IC Socket"
$ws.Range("B4").Value = "ICS28N"
$ws.Range("C4").Value = 0.17
$ws.Range("C4").Style = "Currency"
$ws.Range("D4").Value = "Futurlec"
$ws.Range("E4").Value = "Ok"

# Row 8: 3mm LED
$ws.Range("A8").Value = "3mm LED"
$ws.Range("B8").Value = "LED3R"
$ws.Range("C8").Value = 0.08
$ws.Range("C8").Style = "Currency"
$ws.Range("D8").Value = "Futurlec"
$ws.Range("E8").Value = "Ok"

# Row 9: 100R Resistor
$ws.Range("A9").Value = "100R Resistor"
$ws.Range("B9").Value = "R100R14W"
$ws.Range("C9").Formula = "=0.11/10"
$ws.Range("C9").Style = "Currency"
$ws.Range("D9").Value = "Futurlec"
$ws.Range("E9").Value = "Ok"

# ------------------------------------------------------------------
# Fill in the new rows for Table2
# ------------------------------------------------------------------

# Row 24: IC Socket
$ws.Range("A24").Value = "IC Socket"
$ws.Range("B24").Value = "ICS28N"
$ws.Range("C24").Value = 0.17
$ws.Range("C24").Style = "Currency"
$ws.Range("D24").Value = "Futurlec"
$ws.Range("E24").Value = "Ok"

# Row 29: 3mm LED
$ws.Range("A29").Value = "3mm LED"
$ws.Range("B29").Value = "LED3R"
$ws.Range("C29").Value = 0.08
$ws.Range("C29").Style = "Currency"
$ws.Range("D29").Value = "Futurlec"
$ws.Range("E29").Value = "Ok"

# Row 30: 100R Resistor
$ws.Range("A30").Value = "100R Resistor"
$ws.Range("B30").Value = "R100R14W"
$ws.Range("C30").Formula = "=0.11/10"
$ws.Range("C30").Style = "Currency"
$ws.Range("D30").Value = "Futurlec"
$ws.Range("E30").Value = "Ok"

# ------------------------------------------------------------------
# Resize the two tables (ListObjects) so their ref / autoFilter ranges
# and totals row cover the newly inserted rows.
# ------------------------------------------------------------------
$table1 = $ws.ListObjects.Item("Table1")
$table1.Resize($ws.Range("A2:E19"))

$table2 = $ws.ListObjects.Item("Table2")
$table2.Resize($ws.Range("A22:E33"))

Write-Host "Done updating master parts list"
